$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.153.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.573.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.567.50"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.135.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.568.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.976.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("B23").Value = "RenderToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.75%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "658.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.115"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.413"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.04%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +7.31%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.92%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.146.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.66%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.134"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.131"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.87%  "
